$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42641.54246527778
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9964.52
$ws.Range("D4").Value = 9971
$ws.Range("E4").Value = 79.319999999999993
$ws.Range("F4").Value = 79.22
$ws.Range("G4").Value = $false
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"
$ws.Range("H4").Value = -0.13
$ws.Range("I4").Value = $false
